$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2036.5
$ws.Range("I28").Value = 1312.5454
$ws.Range("J28").Value = 10000
$ws.Range("K28").Value = 1312.5454
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = -827.5454
$ws.Range("N28").Value = -10970
$ws.Range("H74").Value = 2429034.5
$ws.Range("I74").Value = 3399028.5
$ws.Range("J74").Value = 4050
$ws.Range("K74").Value = 3399028.5
$ws.Range("L74").Value = 4050
$ws.Range("M74").Value = -3398092.5
$ws.Range("N74").Value = -5922
$ws.Range("H77").Value = 2429034.5
$ws.Range("I77").Value = 3399028.5
$ws.Range("J77").Value = 4050
$ws.Range("K77").Value = 16995142.5
$ws.Range("L77").Value = 20250
$ws.Range("M77").Value = -16990462.5
$ws.Range("N77").Value = -29610
$ws.Range("H81").Value = 39800
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 39800
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 39800
$ws.Range("N81").Value = -41796
$ws.Range("H84").Value = 39800
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 39800
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 119400
$ws.Range("N84").Value = -129384
$ws.Range("H98").Value = 1515.7858
$ws.Range("I98").Value = 1589.6923
$ws.Range("J98").Value = 555
$ws.Range("K98").Value = 1589.6923
$ws.Range("L98").Value = 555
$ws.Range("M98").Value = -91.69229999999993
$ws.Range("N98").Value = -3551
$ws.Range("H122").Value = 1515.7858
$ws.Range("I122").Value = 1589.6923
$ws.Range("J122").Value = 555
$ws.Range("K122").Value = 4769.0769
$ws.Range("L122").Value = 1665
$ws.Range("M122").Value = -2319.0769
$ws.Range("N122").Value = -6565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2288.0688
$ws.Range("I132").Value = 1650.1428
$ws.Range("J132").Value = 2883.4666
$ws.Range("K132").Value = 4950.428400000001
$ws.Range("L132").Value = 8650.399800000001
$ws.Range("M132").Value = -2420.428400000001
$ws.Range("N132").Value = -13710.3998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 35716556
$ws.Range("I86").Value = 43480280
$ws.Range("J86").Value = 3439.8
$ws.Range("K86").Value = 43480280
$ws.Range("L86").Value = 3439.8
$ws.Range("M86").Value = -43479157
$ws.Range("N86").Value = -5685.8
$ws.Range("H89").Value = 35716556
$ws.Range("I89").Value = 43480280
$ws.Range("J89").Value = 3439.8
$ws.Range("K89").Value = 217401400
$ws.Range("L89").Value = 17199
$ws.Range("M89").Value = -217395784
$ws.Range("N89").Value = -28431

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7522.684
$ws.Range("I86").Value = 6773
$ws.Range("J86").Value = 8553.5
$ws.Range("K86").Value = 6773
$ws.Range("L86").Value = 8553.5
$ws.Range("M86").Value = -5650
$ws.Range("N86").Value = -10799.5
$ws.Range("H89").Value = 7522.684
$ws.Range("I89").Value = 6773
$ws.Range("J89").Value = 8553.5
$ws.Range("K89").Value = 33865
$ws.Range("L89").Value = 42767.5
$ws.Range("M89").Value = -28249
$ws.Range("N89").Value = -53999.5
$ws.Range("H127").Value = 26300
$ws.Range("I127").Value = 7800
$ws.Range("J127").Value = 44800
$ws.Range("K127").Value = 7800
$ws.Range("L127").Value = 44800
$ws.Range("M127").Value = -2840
$ws.Range("N127").Value = -54720
$ws.Range("H134").Value = 38462910
$ws.Range("I134").Value = 1381.7778
$ws.Range("J134").Value = 125001350
$ws.Range("K134").Value = 4145.3334
$ws.Range("L134").Value = 375004050
$ws.Range("M134").Value = -1610.3334
$ws.Range("N134").Value = -375009120

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 53875.26
$ws.Range("I139").Value = 53875.26
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 161625.78
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -156485.78
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 11725
$ws.Range("I46").Value = 5025
$ws.Range("J46").Value = 18425
$ws.Range("K46").Value = 5025
$ws.Range("L46").Value = 18425
$ws.Range("M46").Value = -4869
$ws.Range("N46").Value = -18737
$ws.Range("H57").Value = 15013.75
$ws.Range("J57").Value = 20500
$ws.Range("L57").Value = 20500
$ws.Range("N57").Value = -22140
$ws.Range("H70").Value = 4215.606
$ws.Range("I70").Value = 4346.7915
$ws.Range("J70").Value = 3865.7778
$ws.Range("K70").Value = 4346.7915
$ws.Range("L70").Value = 3865.7778
$ws.Range("M70").Value = -4076.7915
$ws.Range("N70").Value = -4405.7778
$ws.Range("H73").Value = 4215.606
$ws.Range("I73").Value = 4346.7915
$ws.Range("J73").Value = 3865.7778
$ws.Range("K73").Value = 4346.7915
$ws.Range("L73").Value = 3865.7778
$ws.Range("M73").Value = -3410.7915
$ws.Range("N73").Value = -5737.7778
$ws.Range("H119").Value = 24750
$ws.Range("J119").Value = 24750
$ws.Range("L119").Value = 24750
$ws.Range("N119").Value = -34426

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1316.6774
$ws.Range("I16").Value = 1280.2106
$ws.Range("J16").Value = 1374.4166
$ws.Range("K16").Value = 1280.2106
$ws.Range("L16").Value = 1374.4166
$ws.Range("M16").Value = -1110.2106
$ws.Range("N16").Value = -1714.4166
$ws.Range("H22").Value = 670.2
$ws.Range("J22").Value = 850.3333
$ws.Range("L22").Value = 850.3333
$ws.Range("N22").Value = -1440.3333
$ws.Range("H27").Value = 670.2
$ws.Range("J27").Value = 850.3333
$ws.Range("L27").Value = 850.3333
$ws.Range("N27").Value = -1064.3333
$ws.Range("H40").Value = 919947.06
$ws.Range("I40").Value = 1264008.6
$ws.Range("J40").Value = 2449.6667
$ws.Range("K40").Value = 1264008.6
$ws.Range("L40").Value = 2449.6667
$ws.Range("M40").Value = -1263872.6
$ws.Range("N40").Value = -2721.6667
$ws.Range("H46").Value = 1391.091
$ws.Range("J46").Value = 1067
$ws.Range("L46").Value = 1067
$ws.Range("N46").Value = -1443

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 80007
$ws.Range("J15").Value = 80007
$ws.Range("L15").Value = 80007
$ws.Range("N15").Value = -80583
$ws.Range("H54").Value = 6690
$ws.Range("I54").Value = 6690
$ws.Range("K54").Value = 6690
$ws.Range("M54").Value = -6170
$ws.Range("H122").Value = 1070.5714
$ws.Range("I122").Value = 1024
$ws.Range("K122").Value = 3072
$ws.Range("M122").Value = -622
